$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9 - shifts existing rows 9-15 down to 10-16
$ws.Rows.Item(9).Insert()

# Copy formatting from row 8 (identical border/style pattern) into the new row 9
$ws.Range("A8:G8").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 9 with the Plot2DHists_ShearWavy entry
$ws.Range("B9").Value = "Plot2DHists_ShearWavy"
$ws.Range("C9").Value = "Create Gaussian decomposition histograms"
$ws.Range("D9").Value = "ogdenShearVSI-data\Decomposition_Sensititivity\22-1212-Shear_Wavy\sensitivity.mat"
$ws.Range("E9").Value = "Same as param_decoup_[main/nopar] output"
$ws.Range("F9").Value = "ogdenShearVSI-data\Decomposition_Sensititivity\22-12-Shear_Wavy (Folder)"
$ws.Range("G9").Value = "2D histograms of k vs lam for wavy shear simulation"

# Row height for new row
$ws.Rows.Item(9).RowHeight = 45

# Update frozen pane / selection view
$ws.Range("G9").Select() | Out-Null
